$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 352.91666
$ws.Range("I33").Value2 = 218.5
$ws.Range("K33").Value2 = 218.5
$ws.Range("M33").Value2 = 10.5
$ws.Range("H70").Value2 = 1893.75
$ws.Range("I70").Value2 = 1789
$ws.Range("J70").Value2 = 1908.7142
$ws.Range("K70").Value2 = 5367
$ws.Range("L70").Value2 = 5726.142599999999
$ws.Range("M70").Value2 = -5097
$ws.Range("N70").Value2 = -6266.142599999999
$ws.Range("H73").Value2 = 1893.75
$ws.Range("I73").Value2 = 1789
$ws.Range("J73").Value2 = 1908.7142
$ws.Range("K73").Value2 = 5367
$ws.Range("L73").Value2 = 5726.142599999999
$ws.Range("M73").Value2 = -4431
$ws.Range("N73").Value2 = -7598.142599999999
$ws.Range("H127").Value2 = 3311.4736
$ws.Range("I127").Value2 = 3480.4707
$ws.Range("J127").Value2 = 1875
$ws.Range("K127").Value2 = 10441.4121
$ws.Range("L127").Value2 = 5625
$ws.Range("M127").Value2 = -5481.4121
$ws.Range("N127").Value2 = -15545
$ws.Range("H129").Value2 = 465857.88
$ws.Range("I129").Value2 = 537124.4399999999
$ws.Range("J129").Value2 = 2625
$ws.Range("K129").Value2 = 1611373.32
$ws.Range("L129").Value2 = 7875
$ws.Range("M129").Value2 = -1606373.32
$ws.Range("N129").Value2 = -17875
$ws.Range("H132").Value2 = 6644.2925
$ws.Range("I132").Value2 = 6937.2896
$ws.Range("K132").Value2 = 20811.8688
$ws.Range("M132").Value2 = -18281.8688
$ws.Range("H137").Value2 = 5223.7407
$ws.Range("I137").Value2 = 1299.6052
$ws.Range("K137").Value2 = 3898.8156
$ws.Range("M137").Value2 = -1348.8156
$ws.Range("H138").Value2 = 1979.0834
$ws.Range("J138").Value2 = 2054.261
$ws.Range("L138").Value2 = 6162.782999999999
$ws.Range("N138").Value2 = -16442.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 6204.364
$ws.Range("I61").Value2 = 7155.5
$ws.Range("J61").Value2 = 3668
$ws.Range("K61").Value2 = 7155.5
$ws.Range("L61").Value2 = 3668
$ws.Range("M61").Value2 = -6943.5
$ws.Range("N61").Value2 = -4092
$ws.Range("H74").Value2 = 16220.929
$ws.Range("I74").Value2 = 2098.182
$ws.Range("J74").Value2 = 68004.336
$ws.Range("K74").Value2 = 2098.182
$ws.Range("L74").Value2 = 68004.336
$ws.Range("M74").Value2 = -1224.182
$ws.Range("N74").Value2 = -69752.336
$ws.Range("H77").Value2 = 16220.929
$ws.Range("I77").Value2 = 2098.182
$ws.Range("J77").Value2 = 68004.336
$ws.Range("K77").Value2 = 10490.91
$ws.Range("L77").Value2 = 340021.68
$ws.Range("M77").Value2 = -6122.91
$ws.Range("N77").Value2 = -348757.68
$ws.Range("H97").Value2 = 2058.5715
$ws.Range("I97").Value2 = 1538.1818
$ws.Range("K97").Value2 = 1538.1818
$ws.Range("M97").Value2 = -1042.1818
$ws.Range("H132").Value2 = 1471854.9
$ws.Range("I132").Value2 = 1563782.5
$ws.Range("K132").Value2 = 4691347.5
$ws.Range("M132").Value2 = -4688817.5
$ws.Range("H133").Value2 = 74744.75
$ws.Range("J133").Value2 = 74744.75
$ws.Range("L133").Value2 = 74744.75
$ws.Range("N133").Value2 = -79804.75
$ws.Range("H136").Value2 = 6204.364
$ws.Range("I136").Value2 = 7155.5
$ws.Range("J136").Value2 = 3668
$ws.Range("K136").Value2 = 21466.5
$ws.Range("L136").Value2 = 11004
$ws.Range("M136").Value2 = -18916.5
$ws.Range("N136").Value2 = -16104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 1959.1875
$ws.Range("I86").Value2 = 1677.8334
$ws.Range("J86").Value2 = 2803.25
$ws.Range("K86").Value2 = 1677.8334
$ws.Range("L86").Value2 = 2803.25
$ws.Range("M86").Value2 = -554.8334
$ws.Range("N86").Value2 = -5049.25
$ws.Range("H89").Value2 = 1959.1875
$ws.Range("I89").Value2 = 1677.8334
$ws.Range("J89").Value2 = 2803.25
$ws.Range("K89").Value2 = 8389.166999999999
$ws.Range("L89").Value2 = 14016.25
$ws.Range("M89").Value2 = -2773.166999999999
$ws.Range("N89").Value2 = -25248.25
$ws.Range("H94").Value2 = 3734.75
$ws.Range("I94").Value2 = 2980.8
$ws.Range("K94").Value2 = 2980.8
$ws.Range("M94").Value2 = -2529.8
$ws.Range("H134").Value2 = 10800.066
$ws.Range("I134").Value2 = 3818.2727
$ws.Range("J134").Value2 = 30000
$ws.Range("K134").Value2 = 11454.8181
$ws.Range("L134").Value2 = 90000
$ws.Range("M134").Value2 = -8919.8181
$ws.Range("N134").Value2 = -95070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value2 = 54285.715
$ws.Range("J53").Value2 = 50000
$ws.Range("L53").Value2 = 50000
$ws.Range("N53").Value2 = -51214
$ws.Range("H58").Value2 = 11706.786
$ws.Range("J58").Value2 = 18239.428
$ws.Range("L58").Value2 = 18239.428
$ws.Range("N58").Value2 = -18645.428
$ws.Range("H99").Value2 = 4145335.8
$ws.Range("I99").Value2 = 29108.5
$ws.Range("K99").Value2 = 29108.5
$ws.Range("M99").Value2 = -27610.5
$ws.Range("H126").Value2 = 4145335.8
$ws.Range("I126").Value2 = 29108.5
$ws.Range("K126").Value2 = 87325.5
$ws.Range("M126").Value2 = -84855.5
$ws.Range("H136").Value2 = 11706.786
$ws.Range("J136").Value2 = 18239.428
$ws.Range("L136").Value2 = 54718.284
$ws.Range("N136").Value2 = -59818.284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1526.8182
$ws.Range("J5").Value2 = 1661
$ws.Range("L5").Value2 = 4983
$ws.Range("N5").Value2 = -5207
$ws.Range("H81").Value2 = 2082.25
$ws.Range("J81").Value2 = 2398.7
$ws.Range("L81").Value2 = 7196.099999999999
$ws.Range("N81").Value2 = -9442.099999999999
$ws.Range("H84").Value2 = 2082.25
$ws.Range("J84").Value2 = 2398.7
$ws.Range("L84").Value2 = 21588.3
$ws.Range("N84").Value2 = -32820.3
$ws.Range("H92").Value2 = 353.75
$ws.Range("J92").Value2 = 320
$ws.Range("L92").Value2 = 960
$ws.Range("N92").Value2 = -3456
$ws.Range("H132").Value2 = 2060
$ws.Range("I132").Value2 = 1500
$ws.Range("J132").Value2 = 2200
$ws.Range("K132").Value2 = 13500
$ws.Range("L132").Value2 = 19800
$ws.Range("M132").Value2 = -10970
$ws.Range("N132").Value2 = -24860
$ws.Range("H133").Value2 = 7068.6113
$ws.Range("I133").Value2 = 6664.6875
$ws.Range("J133").Value2 = 10300
$ws.Range("K133").Value2 = 19994.0625
$ws.Range("L133").Value2 = 30900
$ws.Range("M133").Value2 = -14934.0625
$ws.Range("N133").Value2 = -41020
$ws.Range("H135").Value2 = 1526.8182
$ws.Range("J135").Value2 = 1661
$ws.Range("L135").Value2 = 14949
$ws.Range("N135").Value2 = -20019

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value2 = 44000
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("H65").Value2 = 44000
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("H80").Value2 = 1394.85
$ws.Range("I80").Value2 = 1099.1818
$ws.Range("J80").Value2 = 1756.2222
$ws.Range("K80").Value2 = 1099.1818
$ws.Range("L80").Value2 = 1756.2222
$ws.Range("M80").Value2 = -101.1818000000001
$ws.Range("N80").Value2 = -3752.2222
$ws.Range("H83").Value2 = 1394.85
$ws.Range("I83").Value2 = 1099.1818
$ws.Range("J83").Value2 = 1756.2222
$ws.Range("K83").Value2 = 5495.909000000001
$ws.Range("L83").Value2 = 8781.110999999999
$ws.Range("M83").Value2 = -503.9090000000006
$ws.Range("N83").Value2 = -18765.111
$ws.Range("H122").Value2 = 3455.9565
$ws.Range("I122").Value2 = 2906.2
$ws.Range("K122").Value2 = 8718.599999999999
$ws.Range("M122").Value2 = -6268.599999999999
$ws.Range("H123").Value2 = 49526
$ws.Range("J123").Value2 = 49526
$ws.Range("L123").Value2 = 49526
$ws.Range("N123").Value2 = -54426
$ws.Range("H126").Value2 = 6194.1177
$ws.Range("J126").Value2 = 2966.6667
$ws.Range("L126").Value2 = 8900.000100000001
$ws.Range("N126").Value2 = -13840.0001
$ws.Range("H140").Value2 = 79999
$ws.Range("J140").Value2 = 79999
$ws.Range("L140").Value2 = 79999
$ws.Range("N140").Value2 = -90359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 913.0833
$ws.Range("I16").Value2 = 677.9091
$ws.Range("K16").Value2 = 677.9091
$ws.Range("M16").Value2 = -507.9091
$ws.Range("H122").Value2 = 4121.4165
$ws.Range("I122").Value2 = 6000
$ws.Range("K122").Value2 = 18000
$ws.Range("M122").Value2 = -15550
$ws.Range("H132").Value2 = 8336583.5
$ws.Range("I132").Value2 = 50000000
$ws.Range("J132").Value2 = 3900
$ws.Range("K132").Value2 = 150000000
$ws.Range("L132").Value2 = 11700
$ws.Range("M132").Value2 = -149997470
$ws.Range("N132").Value2 = -16760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value2 = 350.2857
$ws.Range("I4").Value2 = 350.2857
$ws.Range("K4").Value2 = 350.2857
$ws.Range("M4").Value2 = -237.2857
$ws.Range("H62").Value2 = 9998
$ws.Range("I62").Value2 = 9995
$ws.Range("K62").Value2 = 9995
$ws.Range("M62").Value2 = -9371
$ws.Range("H65").Value2 = 9998
$ws.Range("I65").Value2 = 9995
$ws.Range("K65").Value2 = 49975
$ws.Range("M65").Value2 = -46855
$ws.Range("H81").Value2 = 2427.5789
$ws.Range("I81").Value2 = 1882.75
$ws.Range("K81").Value2 = 3765.5
$ws.Range("M81").Value2 = -2704.5
$ws.Range("H84").Value2 = 2427.5789
$ws.Range("I84").Value2 = 1882.75
$ws.Range("K84").Value2 = 18827.5
$ws.Range("M84").Value2 = -13523.5
$ws.Range("H123").Value2 = 53332.332
$ws.Range("J123").Value2 = 53332.332
$ws.Range("L123").Value2 = 53332.332
$ws.Range("N123").Value2 = -63132.332
$ws.Range("H129").Value2 = 39999.5
$ws.Range("J129").Value2 = 39999.5
$ws.Range("L129").Value2 = 39999.5
$ws.Range("N129").Value2 = -49999.5
